$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds a "wide" layout: one label per row in column A
# (Brown/Green/Blue) followed by that series' numeric readings running
# across the row (B.. onward, variable length per row).
#
#   A1=Brown  B1..I1 = 8 values
#   A2=Green  B2..F2 = 5 values
#   A3=Blue   B3..G3 = 6 values
#
# Re-shape this into a "tall" / tidy layout: the three labels become column
# headers in row 1 (A1:C1) and each series' values run down its own column.

$labels = @()
$series = @()
for ($r = 1; $r -le 3; $r++) {
    $labels += $ws.Cells.Item($r, 1).Value2
    $vals = @()
    $c = 2
    while ($true) {
        $v = $ws.Cells.Item($r, $c).Value2
        if ($v -eq $null -or $v -eq "") { break }
        $vals += $v
        $c += 1
    }
    $series += ,$vals
}

# Wipe the old layout before writing the new one.
$ws.Cells.Clear()

# Row 1: headers (still text, now shared-string refs like before).
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $labels[$i]
}

# Rows 2..: each former row becomes a column of values.
for ($col = 0; $col -lt $series.Length; $col++) {
    $vals = $series[$col]
    for ($row = 0; $row -lt $vals.Length; $row++) {
        $ws.Cells.Item($row + 2, $col + 1).Value = $vals[$row]
    }
}

$ws.Range("A1:C9").Select()
